# Apply the weekly cryptos-list refresh (GitHub Actions scraper run).
# Updates Price (D) and Volume(1h) (E) columns for the changed rows, and
# also restores rows 26/27 (Dai <-> Hedera got reordered upstream) with
# their new Coin/Link/Price/Volume values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.258.15'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').Value = '3.214.21'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''606.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').Value = '''155.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.53%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.212.58'
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('E9').Value = '  -2.36%  '
$ws.Range('D10').Value = '''0.160'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('E11').Value = '  -3.91%  '
$ws.Range('D12').Value = '''0.500'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.63%  '
$ws.Range('D13').Value = '''0.0000267'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.38%  '
$ws.Range('D14').Value = '''38.23'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.17%  '
$ws.Range('D15').Value = '3.737.78'
$ws.Range('E15').Value = '  -0.19%  '
$ws.Range('D16').Value = '66.345.74'
$ws.Range('E16').Value = '  -0.95%  '
$ws.Range('D17').Value = '3.211.15'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '''7.24'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.70%  '
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').Value = '''505.80'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.78%  '
$ws.Range('D21').Value = '''15.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.36%  '
$ws.Range('D22').Value = '''0.727'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.78%  '
$ws.Range('D23').Value = '''7.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.16%  '
$ws.Range('D24').Value = '''14.49'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.26%  '
$ws.Range('D25').Value = '''84.89'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '''1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('B27').Value = 'Hedera'
$ws.Range('C27').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D27').Value = '''0.148'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +62.37%  '
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('D29').Value = '''8.97'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.04%  '
$ws.Range('D30').Value = '''2.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.03%  '
$ws.Range('D31').Value = '''6.91'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.94%  '
$ws.Range('E32').Value = '  -5.00%  '
$ws.Range('D33').Value = '''28.23'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.84%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').Value = '''1.17'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.48%  '
$ws.Range('D36').Value = '''6.37'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.78%  '
$ws.Range('D37').Value = '''55.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('D38').Value = '''497.74'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.44%  '
$ws.Range('D39').Value = '0.0₃0768'
$ws.Range('E39').Value = '  +10.63%  '
$ws.Range('D40').Value = '''0.0417'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.74%  '
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('D42').Value = '''3.01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.46%  '
$ws.Range('D43').Value = '''8.70'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.66%  '
$ws.Range('E44').Value = '  -3.48%  '
$ws.Range('D45').Value = '2.921.24'
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('E46').Value = '  -1.86%  '
$ws.Range('D47').Value = '''27.91'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.11%  '
$ws.Range('D48').Value = '''2.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.87%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  -1.06%  '
$ws.Range('D51').Value = '''121.53'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.38%  '
